$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date update
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicated "Contact" row (old row 11), shifting rows 12-22 up by one
$ws.Rows.Item(11).Delete()

# Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# New Jurisdiction row (replaces the old duplicate Contact row content)
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive value -- must be stored as literal text "true", not boolean TRUE.
# A direct .Value = "true" assignment gets auto-typed to a Boolean by Excel, so
# build it via a formula and then flatten the formula to a static value with a
# copy / paste-values round trip (keeps the string type and original style).
$caseSensitiveCell = $ws.Range("B14")
$caseSensitiveCell.Formula = '="tru"&"e"'
$caseSensitiveCell.Copy()
$caseSensitiveCell.PasteSpecial(-4163)
